$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes the current rows 2-14 down to rows 3-15)
# and populate it with the new top candidate record (Candidate ID 1).
$ws.Rows.Item(2).Insert()

$ws.Range("B2").Value = 1
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2024-01-02"
$ws.Range("D2").Value = "java"
$ws.Range("E2").Value = "surmakkkk"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "55667788"
$ws.Range("G2").Value = "surmkajjjj@gmail.com"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "16"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "15"
$ws.Range("J2").Value = "umbrala corporation"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "3"
$ws.Range("L2").Value = "upgraded for further interview level"
$ws.Range("M2").Value = "shortlisted"

# Append new row 16 - Candidate ID 321
$ws.Range("B16").Value = 321
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "2024-01-02"
$ws.Range("D16").Value = "java"
$ws.Range("E16").Value = "surmakkkk"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "55667788"
$ws.Range("G16").Value = "surmkajjjj@gmail.com"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "16"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "15"
$ws.Range("J16").Value = "umbrala corporation"
$ws.Range("K16").NumberFormat = "@"
$ws.Range("K16").Value = "3"
$ws.Range("L16").Value = "upgraded for further interview level"
$ws.Range("M16").Value = "shortlisted"

# Append new row 17 - Candidate ID 322
$ws.Range("B17").Value = 322
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "2024-01-02"
$ws.Range("D17").Value = "java"
$ws.Range("E17").Value = "surmakkkk"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "55667788"
$ws.Range("G17").Value = "surmkajjjj@gmail.com"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "16"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "15"
$ws.Range("J17").Value = "umbrala corporation"
$ws.Range("K17").NumberFormat = "@"
$ws.Range("K17").Value = "3"
$ws.Range("L17").Value = "upgraded for further interview level"
$ws.Range("M17").Value = "shortlisted"

# Append new row 18 - Candidate ID 323
$ws.Range("B18").Value = 323
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "2024-01-02"
$ws.Range("D18").Value = "java"
$ws.Range("E18").Value = "updated name"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "55667788"
$ws.Range("G18").Value = "surmkaj@gmail.com"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "16"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "15"
$ws.Range("J18").Value = "umbrala corporation"
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value = "3"
$ws.Range("L18").Value = "upgraded for further interview level"
$ws.Range("M18").Value = "shortlisted"
